# feat: add 2022-Q4 data
#
# Inserts a new quarterly snapshot sheet "2022-Q4" (positioned right after the
# "总计" summary sheet, i.e. before "2022-Q3") and records it as a new top row
# in the "总计" summary sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Update the "总计" summary sheet: insert a new row for 2022-Q4 ---
$summary = $wb.Worksheets.Item("总计")

# Push existing data rows (2022-Q3 / 2022-Q2 / 2021-Q2) down by one row.
$summary.Rows.Item(2).Insert()

# The freshly-inserted row borrows formatting from the header row above it;
# strip that back to plain (matching the look of the other data rows), then
# restore column A's "index" style by copying it from the row beneath.
$summary.Range("B2:D2").ClearFormats()
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# Fill in the new 2022-Q4 row and renumber the index column (A) sequentially.
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 2.39

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3

# --- 2. Create the new "2022-Q4" detail sheet ---
# Duplicate the existing "2022-Q3" sheet (so formatting/column layout match)
# and place the copy immediately before it, then rename + rewrite its values.
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Row 2 - 008763 / 天弘越南市场股票（QDII）A
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "21.47"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "92.63"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "6.34"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "1.3612"
$q4.Range("H2").Value = 6

# Row 3 - 008764 / 天弘越南市场股票（QDII）C
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "16.17"
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "92.63"
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "6.34"
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "1.0252"
$q4.Range("H3").Value = 6
